# daily auto push: 2025-10-09 02:00 UTC
# Append a new daily record row to the sheet's data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 83

# Column A holds the date as literal text (matches the existing rows, which
# store "YYYY/MM/DD" as plain text rather than a real date value). A leading
# apostrophe forces text entry instead of Excel's automatic date parsing;
# resetting the style afterward keeps the cell on the default (unstyled)
# format like every other data row.
$ws.Cells.Item($newRow, 1).Formula = "'2025/10/09"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "木"
$ws.Cells.Item($newRow, 3).Value = 11
$ws.Cells.Item($newRow, 4).Value = 141
